$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.891504666666667
$ws.Range("H2").Value = 8.674514
$ws.Range("I2").Value = 0.1213590456377548
$ws.Range("J2").Value = 0.1213590456377548
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 147.4213356666667
$ws.Range("N2").Value = 442.264007
$ws.Range("O2").Value = 0.9507885170992249
$ws.Range("P2").Value = 0.950788517099225
$ws.Range("Q2").Value = 426.2694800463997
$ws.Range("R2").Value = 3836.425320417598
$ws.Range("S2").Value = 0.115386787038498
$ws.Range("T2").Value = 0.115386787038498

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.891504666666667
$ws.Range("H3").Value = 8.674514
$ws.Range("I3").Value = 0.1213590456377548
$ws.Range("J3").Value = 0.1213590456377548
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.340788333333334
$ws.Range("N3").Value = 7.022365000000001
$ws.Range("O3").Value = 0.01509682881537204
$ws.Range("P3").Value = 0.01509682881537204
$ws.Range("Q3").Value = 6.768400389512223
$ws.Range("R3").Value = 60.91560350561001
$ws.Range("S3").Value = 0.001832136737190106
$ws.Range("T3").Value = 0.001832136737190106

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.891504666666667
$ws.Range("H4").Value = 8.674514
$ws.Range("I4").Value = 0.1213590456377548
$ws.Range("J4").Value = 0.1213590456377548
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.289533666666667
$ws.Range("N4").Value = 15.868601
$ws.Range("O4").Value = 0.03411465408540306
$ws.Range("P4").Value = 0.03411465408540307
$ws.Range("Q4").Value = 15.29471128165711
$ws.Range("R4").Value = 137.652401534914
$ws.Range("S4").Value = 0.004140121862066647
$ws.Range("T4").Value = 0.004140121862066648

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 12.04042966666667
$ws.Range("H5").Value = 36.121289
$ws.Range("I5").Value = 0.505347637947847
$ws.Range("J5").Value = 0.505347637947847
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 147.4213356666667
$ws.Range("N5").Value = 442.264007
$ws.Range("O5").Value = 0.9507885170992249
$ws.Range("P5").Value = 0.950788517099225
$ws.Range("Q5").Value = 1775.016223460558
$ws.Range("R5").Value = 15975.14601114503
$ws.Range("S5").Value = 0.4804787313040294
$ws.Range("T5").Value = 0.4804787313040295

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.04042966666667
$ws.Range("H6").Value = 36.121289
$ws.Range("I6").Value = 0.505347637947847
$ws.Range("J6").Value = 0.505347637947847
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.340788333333334
$ws.Range("N6").Value = 7.022365000000001
$ws.Range("O6").Value = 0.01509682881537204
$ws.Range("P6").Value = 0.01509682881537204
$ws.Range("Q6").Value = 28.1840972920539
$ws.Range("R6").Value = 253.6568756284851
$ws.Range("S6").Value = 0.00762914678235125
$ws.Range("T6").Value = 0.00762914678235125

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.04042966666667
$ws.Range("H7").Value = 36.121289
$ws.Range("I7").Value = 0.505347637947847
$ws.Range("J7").Value = 0.505347637947847
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.289533666666667
$ws.Range("N7").Value = 15.868601
$ws.Range("O7").Value = 0.03411465408540306
$ws.Range("P7").Value = 0.03411465408540307
$ws.Range("Q7").Value = 63.68825808296546
$ws.Range("R7").Value = 573.1943227466892
$ws.Range("S7").Value = 0.0172397598614663
$ws.Range("T7").Value = 0.01723975986146631

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.894099
$ws.Range("H8").Value = 26.682297
$ws.Range("I8").Value = 0.3732933164143983
$ws.Range("J8").Value = 0.3732933164143982
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 147.4213356666667
$ws.Range("N8").Value = 442.264007
$ws.Range("O8").Value = 0.9507885170992249
$ws.Range("P8").Value = 0.950788517099225
$ws.Range("Q8").Value = 1311.179954131565
$ws.Range("R8").Value = 11800.61958718408
$ws.Range("S8").Value = 0.3549229987566975
$ws.Range("T8").Value = 0.3549229987566975

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.894099
$ws.Range("H9").Value = 26.682297
$ws.Range("I9").Value = 0.3732933164143983
$ws.Range("J9").Value = 0.3732933164143982
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.340788333333334
$ws.Range("N9").Value = 7.022365000000001
$ws.Range("O9").Value = 0.01509682881537204
$ws.Range("P9").Value = 0.01509682881537204
$ws.Range("Q9").Value = 20.81920317471167
$ws.Range("R9").Value = 187.372828572405
$ws.Range("S9").Value = 0.005635545295830679
$ws.Range("T9").Value = 0.005635545295830678

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.894099
$ws.Range("H10").Value = 26.682297
$ws.Range("I10").Value = 0.3732933164143983
$ws.Range("J10").Value = 0.3732933164143982
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.289533666666667
$ws.Range("N10").Value = 15.868601
$ws.Range("O10").Value = 0.03411465408540306
$ws.Range("P10").Value = 0.03411465408540307
$ws.Range("Q10").Value = 47.04563609516634
$ws.Range("R10").Value = 423.4107248564971
$ws.Range("S10").Value = 0.01273477236187011
$ws.Range("T10").Value = 0.01273477236187011

